$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 92, shifting existing rows 92-200 down to 93-201,
# copying formatting from the row above (so the date style on column D is preserved).
$ws.Rows.Item(92).Insert(-4121, 0)  # xlShiftDown = -4121, xlFormatFromLeftOrAbove = 0

# Populate the newly inserted row 92 with the new weekly data record.
$ws.Cells.Item(92, 1).Value = 4
$ws.Cells.Item(92, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(92, 3).Value = "Los Lagos"
$ws.Cells.Item(92, 4).Value = 44579
$ws.Cells.Item(92, 5).Value = 10
$ws.Cells.Item(92, 6).Value = 100112021
$ws.Cells.Item(92, 7).Value = "Ají"
$ws.Cells.Item(92, 8).Value = "Inferno"
$ws.Cells.Item(92, 9).Value = "Primera"
$ws.Cells.Item(92, 10).Value = 80
$ws.Cells.Item(92, 11).Value = 25000
$ws.Cells.Item(92, 12).Value = 25000
$ws.Cells.Item(92, 13).Value = 25000
$ws.Cells.Item(92, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(92, 15).Value = "Región Metropolitana"
$ws.Cells.Item(92, 16).Value = 1667
$ws.Cells.Item(92, 17).Value = 15
$ws.Cells.Item(92, 18).Value = "Hortaliza"
